# Apply the "add example values to sequencing related templates" edit to the
# RNASeq_computational_analysis template sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rnaseq_computational_analysis")

# 1. Rename the last header / table column:
#    "Output [Derived Data File]" -> "Output [Data]"
$ws.Range("AD1").Value = "Output [Data]"

# 2. Fill in new example values for the single remaining data row (row 2).
$ws.Range("F2").Value = "CLC Genomics Workbench"   # Parameter [data filtering software]
$ws.Range("I2").Value = "v11.0.1"                  # Parameter [data filtering software version]
$ws.Range("L2").Value = "quality limit: 0.05"      # Parameter [data filtering software parameters]
$ws.Range("P2").Value = ""                          # Term Source REF (DPBO:0000002) - cleared
$ws.Range("X2").Value = "B73 RefGen_v4"            # Parameter [genome reference sequence]

# 3. Remove the now-superfluous example rows 3-5, shrinking the table/sheet
#    down to a single example data row (A1:AD2).
$ws.Range("A3:A5").EntireRow.Delete()
